$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list price/volume refresh — update Price (D) and Volume(1h) (E) columns
# D-column cells must remain text (prices like "40.123.48" / "294.17"), so force
# the text number format before assigning, matching the source data formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.123.48'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.225.18'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.17'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.99'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.469'
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.80'
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.88'
$ws.Range("E11").Value = '  +6.17%  '
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("E13").Value = '  +3.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.45'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.585.46'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.87'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.270.59'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.739'
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.082.55'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.32'
$ws.Range("E21").Value = '  -5.27%  '
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.74'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.40'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.23'
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.33'
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("E30").Value = '  -6.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.71'
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.89'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  +6.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0717'
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("E39").Value = '  +2.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0995'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.57'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.086.38'
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.74'
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.01'
$ws.Range("E44").Value = '  +7.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.13'
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  +2.66%  '
$ws.Range("E48").Value = '  -11.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.453.69'
$ws.Range("E49").Value = '  +1.06%  '
$ws.Range("E50").Value = '  +2.56%  '
$ws.Range("E51").Value = '  +3.50%  '
